$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "Value_LSTM_CNN" column (C) with the loss statistics ---
$ws.Range("C2").Value2 = 0.000911727768252604
$ws.Range("C3").Value2 = 0.000182569790044063
$ws.Range("C4").Value2 = 0.000636487442534417
$ws.Range("C5").Value2 = 0.00145924591924995
$ws.Range("C6").Value2 = 0.000864890549564734
$ws.Range("C7").Value2 = 0.000636487442534417

# Match column C's number formatting / alignment to column B (0.0000, centered)
$ws.Range("C2:C7").NumberFormat = "0.0000"
$ws.Range("C2:C7").HorizontalAlignment = -4108

# --- Add the "CV" (coefficient of variation) row as formulas ---
$ws.Range("B8").Formula = "=B3/B2"
$ws.Range("C8").Formula = "=C3/C2"
$ws.Range("B8:C8").NumberFormat = "0.0000"
$ws.Range("B8:C8").HorizontalAlignment = -4108

# --- Move the active cell/selection like the authored workbook ---
[void]$ws.Range("C14").Select()
